$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Data for new rows 2-7 (A..H). Columns I and J left blank, as in the diff.
$data = @(
    @(1, 900858550,    "BIBO",    "SOLUTIONS SAS",             7,    "CA", 6456465, 2359386),
    @(1, 900654100,    "CIMAZ",   "S.A.S",                     7,    "CA", 654546,  525870),
    @(1, 1143940723,   "CIMPRE",  "SALUD OCUPACIONAL S.A.S.",  7,    "CA", 5646545, 87451),
    @(3, 31322510,     "DIANA",   "CARINA IMPATA RESTREPO",    1051, "CA", 654654,  84300),
    @(1, 901223156,    "DIGITALTIC", "SAS",                    1057, 46546, 645654, 193970),
    @(3, 7215649,      "DOMINGO", "IGNACIO ROJAS",             1057, "CA", 64654,   97991)
)

$startRow = 2
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $vals = $data[$i]

    $ws.Cells.Item($row, 1).Value = $vals[0]
    $ws.Cells.Item($row, 2).Value = $vals[1]
    $ws.Cells.Item($row, 3).Value = $vals[2]
    $ws.Cells.Item($row, 4).Value = $vals[3]
    $ws.Cells.Item($row, 5).Value = $vals[4]
    $ws.Cells.Item($row, 6).Value = $vals[5]
    $ws.Cells.Item($row, 7).Value = $vals[6]
    $ws.Cells.Item($row, 8).Value = $vals[7]
}
